$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Recorded By" values in column G (rows 2-187), leaving the header
# in G1 and the style of each cell untouched.
$ws.Range("G2:G187").ClearContents()

# Shrink column G from its original width (50) down to 13. Excel's
# ColumnWidth property is offset from the raw OOXML column width by the
# standard ~0.8333 character padding, so subtract that offset to land on
# a stored width of exactly 13.
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666
